$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "food colorant, cosmetic additive, dietary supplement, health food"
$ws.Range("D3").Value = "food colorant, dietary supplement, health food, antioxidant, minerals"
$ws.Range("D5").Value = "food colorant, health food, algae"
$ws.Range("D7").Value = "food colorant, health food, algae, minerals"
$ws.Range("D8").Value = "food colorant, dietary supplement, health food, antioxidant"
$ws.Range("D9").Value = "antioxidant"
$ws.Range("D11").Value = "food colorant, health food"

$wb.Save()
